# InstantiateAgents.xlsx — step-function rewrite
# x/y coordinates are being re-based onto an Earth-centered frame by adding
# the Earth radius (6378.14 km) via live formulas instead of hard literals,
# and the Transporters sheet's "orig" placeholder ("earth") is replaced with
# the actual origin node id for the two rows that used it.

$wb = $excel.ActiveWorkbook

$fix = $wb.Worksheets.Item("FixNodes")
$var = $wb.Worksheets.Item("VarNodes")
$trn = $wb.Worksheets.Item("Transporters")

# ---- FixNodes (FixTest1, row 2): x/y rebased to Earth-centered coords ----
$fix.Range("B2").Formula = "=6378.14+300"
$fix.Range("C2").Formula = "=300+6378.14"

# ---- VarNodes (VarTest1 row2 / VarTest2 row3): x/y rebased ----
$var.Range("B2").Formula = "=200+6378.14"
$var.Range("C2").Formula = "=350+6378.14"
$var.Range("B3").Formula = "=800+6378.14"
$var.Range("C3").Formula = "=950+6378.14"

# ---- Transporters: x/y rebased for all three rows ----
$trn.Range("C2").Formula = "=-(250+6378.14)"
$trn.Range("D2").Formula = "=250+6378.14"
$trn.Range("C3").Formula = "=350+6378.14"
$trn.Range("D3").Formula = "=120+6378.14"
$trn.Range("C4").Formula = "=300+6378.14"
$trn.Range("D4").Formula = "=200+6378.14"

# Transporters "orig" column: replace the placeholder "earth" text with the
# real origin node id for TestTrans1 and TestTrans3 (TestTrans2 already had
# a real node id and is untouched).
$trn.Range("E2").Value = "VarTest1"
$trn.Range("E4").Value = "FixTest1"

# ---- Restore selections on each sheet (also re-marks FixNodes as the
# active tab, since it must be selected last to keep tabSelected="1") ----
$var.Range("C4").Select() | Out-Null
$trn.Range("C26").Select() | Out-Null
$fix.Range("C31").Select() | Out-Null

$wb.Save()
